$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''59.898.77'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '''  +1.92%  '
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = '''2.555.82'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '''  +3.85%  '
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = '''  -0.06%  '
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = '''502.74'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '''  +1.15%  '
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = '''152.03'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '''  -4.48%  '
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Value = '''  +0.58%  '
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = '''0.574'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '''  -6.66%  '
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = '''2.569.30'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '''  +3.18%  '
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = '''6.77'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '''  +7.56%  '
$ws.Range("E10").Style = "Normal"
$ws.Range("E11").Value = '''  -0.81%  '
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = '''0.342'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '''  +1.24%  '
$ws.Range("E12").Style = "Normal"
$ws.Range("E13").Value = '''  +0.58%  '
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = '''3.020.00'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '''  +5.02%  '
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = '''60.006.70'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '''  +2.33%  '
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = '''21.41'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '''  -1.74%  '
$ws.Range("E16").Style = "Normal"
$ws.Range("E17").Value = '''  +1.74%  '
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = '''2.569.33'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '''  +4.13%  '
$ws.Range("E18").Style = "Normal"
$ws.Range("E19").Value = '''  +0.94%  '
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = '''345.04'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '''  +4.33%  '
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = '''10.19'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '''  +0.76%  '
$ws.Range("E21").Style = "Normal"
$ws.Range("E22").Value = '''  +0.56%  '
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = '''1.00'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '''  +0.19%  '
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = '''59.79'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '''  +1.43%  '
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = '''0.416'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '''  +1.22%  '
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = '''0.164'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '''  -1.14%  '
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = '''2.679.43'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '''  +4.35%  '
$ws.Range("E27").Style = "Normal"
$ws.Range("E28").Value = '''  +0.27%  '
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = '''0.0₃0841'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '''  +4.02%  '
$ws.Range("E29").Style = "Normal"
$ws.Range("E30").Value = '''  -0.36%  '
$ws.Range("E30").Style = "Normal"
$ws.Range("E31").Value = '''  +0.31%  '
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = '''155.33'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '''  +2.20%  '
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = '''19.08'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '''  +0.11%  '
$ws.Range("E33").Style = "Normal"
$ws.Range("E34").Value = '''  +0.30%  '
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = '''5.68'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '''  +3.78%  '
$ws.Range("E35").Style = "Normal"
$ws.Range("E36").Value = '''  +2.09%  '
$ws.Range("E36").Style = "Normal"
$ws.Range("E37").Value = '''  +0.35%  '
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = '''0.854'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '''  +20.24%  '
$ws.Range("E38").Style = "Normal"
$ws.Range("B39").Value = '''Fetch.AI'
$ws.Range("B39").Style = "Normal"
$ws.Range("C39").Value = '''https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("C39").Style = "Normal"
$ws.Range("D39").Value = '''0.838'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '''  -1.38%  '
$ws.Range("E39").Style = "Normal"
$ws.Range("B40").Value = '''Filecoin'
$ws.Range("B40").Style = "Normal"
$ws.Range("C40").Value = '''https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("C40").Style = "Normal"
$ws.Range("D40").Value = '''3.75'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '''  +2.25%  '
$ws.Range("E40").Style = "Normal"
$ws.Range("E41").Value = '''  +1.25%  '
$ws.Range("E41").Style = "Normal"
$ws.Range("B42").Value = '''OKB'
$ws.Range("B42").Style = "Normal"
$ws.Range("C42").Value = '''https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("C42").Style = "Normal"
$ws.Range("D42").Value = '''35.41'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '''  +2.71%  '
$ws.Range("E42").Style = "Normal"
$ws.Range("B43").Value = '''Bittensor'
$ws.Range("B43").Style = "Normal"
$ws.Range("C43").Value = '''https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("C43").Style = "Normal"
$ws.Range("D43").Value = '''295.32'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '''  +3.26%  '
$ws.Range("E43").Style = "Normal"
$ws.Range("E44").Value = '''  +2.30%  '
$ws.Range("E44").Style = "Normal"
$ws.Range("E45").Value = '''  -1.95%  '
$ws.Range("E45").Style = "Normal"
$ws.Range("B46").Value = '''FirstDigitalUSD'
$ws.Range("B46").Style = "Normal"
$ws.Range("C46").Value = '''https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("C46").Style = "Normal"
$ws.Range("D46").Value = '''0.997'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '''  +0.60%  '
$ws.Range("E46").Style = "Normal"
$ws.Range("B47").Value = '''Mantle'
$ws.Range("B47").Style = "Normal"
$ws.Range("C47").Value = '''https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("C47").Style = "Normal"
$ws.Range("D47").Value = '''0.612'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '''  +0.35%  '
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = '''19.49'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '''  +7.32%  '
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = '''4.85'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '''  +1.09%  '
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = '''0.0233'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '''  -1.27%  '
$ws.Range("E50").Style = "Normal"
$ws.Range("E51").Value = '''  -0.08%  '
$ws.Range("E51").Style = "Normal"

Write-Host "Done: applied all cell updates"
